$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "0.999") must be
# pre-formatted as Text, otherwise Excel auto-converts the assigned string
# into a numeric value and the original text formatting (trailing zeros,
# "1.00" vs "1", etc.) would be lost.

$ws.Range("D2").Value = "51.426.55"
$ws.Range("E2").Value = "  -15.31%  "
$ws.Range("D3").Value = "2.258.69"
$ws.Range("E3").Value = "  -22.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "431.39"
$ws.Range("E5").Value = "  -18.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.01"
$ws.Range("E6").Value = "  -19.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -17.58%  "
$ws.Range("D9").Value = "2.251.71"
$ws.Range("E9").Value = "  -22.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.12"
$ws.Range("E10").Value = "  -15.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0830"
$ws.Range("E11").Value = "  -22.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.293"
$ws.Range("E12").Value = "  -17.99%  "
$ws.Range("E13").Value = "  -6.90%  "
$ws.Range("D14").Value = "2.624.08"
$ws.Range("E14").Value = "  -23.03%  "
$ws.Range("D15").Value = "51.390.23"
$ws.Range("E15").Value = "  -15.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.15"
$ws.Range("E16").Value = "  -19.49%  "
$ws.Range("E17").Value = "  -20.23%  "
$ws.Range("D18").Value = "2.255.19"
$ws.Range("E18").Value = "  -22.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.83"
$ws.Range("E19").Value = "  -22.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "290.55"
$ws.Range("E20").Value = "  -17.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.38"
$ws.Range("E23").Value = "  -27.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.89"
$ws.Range("E24").Value = "  -24.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "52.48"
$ws.Range("E26").Value = "  -19.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.356"
$ws.Range("E27").Value = "  -20.94%  "
$ws.Range("D28").Value = "2.338.97"
$ws.Range("E28").Value = "  -22.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.133"
$ws.Range("E29").Value = "  -24.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("E31").Value = "  -15.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "142.94"
$ws.Range("E32").Value = "  -6.17%  "
$ws.Range("D33").Value = "0.0₃0618"
$ws.Range("E33").Value = "  -27.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "16.21"
$ws.Range("E34").Value = "  -17.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.28"
$ws.Range("E35").Value = "  -23.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.53"
$ws.Range("E36").Value = "  -18.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.991"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.767"
$ws.Range("E38").Value = "  -22.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value = "  -25.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.949"
$ws.Range("E40").Value = "  -20.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "31.51"
$ws.Range("E41").Value = "  -15.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.12"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.543"
$ws.Range("E43").Value = "  -16.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0483"
$ws.Range("E44").Value = "  -16.89%  "
$ws.Range("E45").Value = "  -19.18%  "
$ws.Range("D46").Value = "1.837.56"
$ws.Range("E46").Value = "  -19.79%  "
$ws.Range("E47").Value = "  -24.45%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0198"
$ws.Range("E48").Value = "  -16.80%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0786"
$ws.Range("E49").Value = "  -14.29%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.87"
$ws.Range("E50").Value = "  -21.56%  "
$ws.Range("B51").Value = "ZEEBU"
$ws.Range("C51").Value = "https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.62"
$ws.Range("E51").Value = "  -5.31%  "
